$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-text labels: replace embedded newlines with literal "<br>" markers ---
$ws.Range("E26:E37").Value = "Supports the GCS if coverage is **Low**<br>Other members: Global South + EU<br>(25-33% of world emissions)"
$ws.Range("E38:E49").Value = "Supports the GCS if coverage is **Mid**<br>Global South + China<br>(56% of world emissions)"
$ws.Range("E50:E61").Value = "Supports the GCS if coverage is **High**<br>Global South + China + EU + various HICs<br>(UK, Japan, Korea, Canada...; 64-72% of emissions)"
$ws.Range("E62:E73").Value = "Supports the GCS if coverage is **High**, **color** variant<br>Global South + China + EU + various HICs<br>+ Distributive effects shown using colors on world map"

# --- Update recomputed statistics (mean / CI_low / CI_high) after removing fast RU respondents ---
$ws.Range("B2").Value = 67.8255122017956
$ws.Range("C2").Value = 66.435391002076
$ws.Range("D2").Value = 69.2156334015152

$ws.Range("B12").Value = 73.9047667329172
$ws.Range("C12").Value = 70.0282519152882
$ws.Range("D12").Value = 77.7812815505462

$ws.Range("B14").Value = 55.3352948193206
$ws.Range("C14").Value = 53.8534640960617
$ws.Range("D14").Value = 56.8171255425794

$ws.Range("B24").Value = 48.6313846276599
$ws.Range("C24").Value = 44.2764656460201
$ws.Range("D24").Value = 52.9863036092996

$ws.Range("B38").Value = 67.1354010141054
$ws.Range("C38").Value = 65.427782749948
$ws.Range("D38").Value = 68.8430192782629

$ws.Range("B48").Value = 63.442019211072
$ws.Range("C48").Value = 57.6976042102314
$ws.Range("D48").Value = 69.1864342119125

$ws.Range("B50").Value = 68.4899159160604
$ws.Range("C50").Value = 66.8593150984853
$ws.Range("D50").Value = 70.1205167336354

$ws.Range("B60").Value = 60.1821932205212
$ws.Range("C60").Value = 54.1620904472528
$ws.Range("D60").Value = 66.2022959937895

$ws.Range("B62").Value = 61.8798692282585
$ws.Range("C62").Value = 60.1485598973454
$ws.Range("D62").Value = 63.6111785591716

$ws.Range("B72").Value = 54.0183622108344
$ws.Range("C72").Value = 47.9013277781348
$ws.Range("D72").Value = 60.135396643534
